$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt - Cilantro".
# It belongs right above the existing row 67, so insert a new row there; this pushes the
# old rows 67..172 down to 68..173 (old row 172 data ends up in the newly created row 173).
$ws.Rows.Item(67).Insert()

# Make sure the date cell in the new row keeps the same date style used by the rest of
# column D (style index 2 / custom date-time numFmt), by copying it from the row below.
$ws.Cells.Item(68, 4).Copy()
$ws.Cells.Item(67, 4).PasteSpecial(-4122) | Out-Null

# Fill in the new record's values
$ws.Cells.Item(67, 1).Value = 4
$ws.Cells.Item(67, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(67, 3).Value = "Los Lagos"
$ws.Cells.Item(67, 4).Value = 44477
$ws.Cells.Item(67, 5).Value = 10
$ws.Cells.Item(67, 6).Value = 100112040
$ws.Cells.Item(67, 7).Value = "Cilantro"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 250
$ws.Cells.Item(67, 11).Value = 10000
$ws.Cells.Item(67, 12).Value = 10000
$ws.Cells.Item(67, 13).Value = 10000
$ws.Cells.Item(67, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(67, 15).Value = "Región Metropolitana"
$ws.Cells.Item(67, 16).Value = 278
$ws.Cells.Item(67, 17).Value = 36
$ws.Cells.Item(67, 18).Value = "Hortaliza"
